# "adding term 2.0 now utf-8"
#
# 1. Bump the Metadata sheet's Version/Date/Contact values.
# 2. Replace the placeholder "I<n>" Value codes on the 12 existing
#    "Include from FSIII*" sheets with their real UUID codes.
# 3. Append 12 new "Include from FSIII 13".."Include from FSIII 24" sheets
#    (clones of the existing template sheet) whose Value column now carries
#    the freed-up "I1".."I12" placeholders.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Metadata sheet updates
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(3, 2).Value  = "2.0.0"                                        # Version
$meta.Cells.Item(8, 2).Value  = "2024-06-03T10:45:43+02:00"                    # Date
$meta.Cells.Item(10, 2).Value = "Kommunernes Landsforening (http://kl.dk)"     # Contact

# ---------------------------------------------------------------------------
# 2. Swap the "I<n>" placeholders on the original 12 sheets for real UUIDs
# ---------------------------------------------------------------------------
$newCodes = @(
    "fa6aa904-d06e-4029-b4c4-13ead04ace27",
    "55670b1e-7a36-46b2-8712-b7536237f22d",
    "25dcedb3-7149-4ef9-a2c3-be30267441fb",
    "1bb534f3-e526-41a9-b9c3-6157ea19c915",
    "3f00a76f-8e7b-4b13-80cc-f2ceef4e51d1",
    "5bfe4bda-2358-41da-946e-1fdaa33d5fe8",
    "01150cdb-6098-48ce-bb61-60967f6bcc37",
    "cc377732-7f14-49b7-8940-1aa07b8884e7",
    "94e9c867-fbc8-4d35-8596-e6b8765b12e8",
    "8c539fd9-7f31-4b4e-8b30-8298c8ab640f",
    "9162d29a-1c7f-4585-8145-8fb4f1a999e3",
    "045fa500-35b0-46b7-97dd-adb60888a8ea"
)

for ($i = 0; $i -lt $newCodes.Length; $i++) {
    # sheets are at workbook index 2..13 ("Include from FSIII", "Include from FSIII 2", ...)
    $ws = $wb.Worksheets.Item($i + 2)
    $ws.Cells.Item(2, 3).Value = $newCodes[$i]
}

# ---------------------------------------------------------------------------
# 3. Add 12 new clone sheets carrying the freed "I1".."I12" placeholders
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(2)   # "Include from FSIII" - same layout/styles for all of them

for ($n = 13; $n -le 24; $n++) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $lastSheet)
    $newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet.Name = "Include from FSIII $n"
    $newSheet.Cells.Item(2, 3).Value = "I" + ($n - 12)
}
